$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.680.40"
$ws.Range("E2").Value = "  +7.52%  "
$ws.Range("D3").Value = "1.745.15"
$ws.Range("E3").Value = "  +5.50%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'334.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.50%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").Value = "'0.3753"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.77%  "
$ws.Range("D8").Value = "'49.24"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.3399"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.07%  "
$ws.Range("D10").Value = "'1.197"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.47%  "
$ws.Range("D11").Value = "'0.07484"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.71%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").Value = "'6.480"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.07%  "
$ws.Range("D14").Value = "'20.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.37%  "
$ws.Range("D15").Value = "'7.133"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +9.04%  "
$ws.Range("D16").Value = "1.746.74"
$ws.Range("E16").Value = "  +5.62%  "
$ws.Range("E17").Value = "  +4.93%  "
$ws.Range("D18").Value = "'0.06703"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.73%  "
$ws.Range("D19").Value = "'83.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.70%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").Value = "'16.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.79%  "
$ws.Range("D22").Value = "'6.197"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.84%  "
$ws.Range("D23").Value = "'13.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.58%  "
$ws.Range("D24").Value = "26.639.95"
$ws.Range("E24").Value = "  +7.55%  "
$ws.Range("D25").Value = "'2.443"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.60%  "
$ws.Range("D26").Value = "'2.482"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.04%  "
$ws.Range("D27").Value = "'1.426"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +20.25%  "
$ws.Range("D28").Value = "'154.48"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.06%  "
$ws.Range("D29").Value = "'19.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.85%  "
$ws.Range("D30").Value = "1.937.74"
$ws.Range("E30").Value = "  +5.41%  "
$ws.Range("D31").Value = "'132.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.94%  "
$ws.Range("D32").Value = "'4.123"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.03%  "
$ws.Range("D33").Value = "'6.133"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.02%  "
$ws.Range("D34").Value = "'0.08670"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.36%  "
$ws.Range("D35").Value = "'1.712"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.60%  "
$ws.Range("D36").Value = "'13.07"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.99%  "
$ws.Range("D37").Value = "'5.450"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.39%  "
$ws.Range("D38").Value = "'0.02365"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.52%  "
$ws.Range("D39").Value = "'0.06324"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.43%  "
$ws.Range("D40").Value = "'0.2187"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.67%  "
$ws.Range("D41").Value = "'8.614"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.65%  "
$ws.Range("D42").Value = "'1.230"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.88%  "
$ws.Range("D43").Value = "'0.6270"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.77%  "
$ws.Range("D44").Value = "'14.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.95%  "
$ws.Range("D45").Value = "'1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("D46").Value = "'3.933"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.56%  "
$ws.Range("D47").Value = "'0.6074"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.38%  "
$ws.Range("D48").Value = "'129.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.24%  "
$ws.Range("D49").Value = "'2.069"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.52%  "
$ws.Range("D50").Value = "'0.07282"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.86%  "
$ws.Range("D51").Value = "'78.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.40%  "
